$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("T2").Value = 0.1426048780487805
$ws.Range("V2").Value = 0.0002220611491829204
$ws.Range("Z2").Value = -0.08228521019851517
$ws.Range("AB2").Value = -370.5520326328386
$ws.Range("AC2").Value = "umolO2/min/m2"
$ws.Range("AD2").Value = -370.5520326328386

$ws.Range("T3").Value = 0.1423414634146342
$ws.Range("V3").Value = 0.0002565454225970831
$ws.Range("Z3").Value = -0.1400753476368244
$ws.Range("AB3").Value = -546.0060297268269
$ws.Range("AC3").Value = "umolO2/min/m2"
$ws.Range("AD3").Value = -546.0060297268269

$ws.Range("T4").Value = 0.1477268292682927
$ws.Range("V4").Value = 0.0002491214197856265
$ws.Range("Z4").Value = -0.1320049077317648
$ws.Range("AB4").Value = -529.8818056085159
$ws.Range("AC4").Value = "umolO2/min/m2"
$ws.Range("AD4").Value = -529.8818056085159

$ws.Range("T5").Value = 0.1470341463414634
$ws.Range("V5").Value = 0.0002082674398172554
$ws.Range("Z5").Value = -0.121577720776101
$ws.Range("AB5").Value = -583.7576957914285
$ws.Range("AC5").Value = "umolO2/min/m2"
$ws.Range("AD5").Value = -583.7576957914285

$ws.Range("T6").Value = 0.1477560975609756
$ws.Range("V6").Value = 0.0001607142857142857
$ws.Range("Z6").Value = -0.1354951947609839
$ws.Range("AB6").Value = -843.0812118461221
$ws.Range("AC6").Value = "umolO2/min/m2"
$ws.Range("AD6").Value = -843.0812118461221

$ws.Range("T7").Value = 0.1455317073170732
$ws.Range("V7").Value = 0.0002247627833421192
$ws.Range("Z7").Value = -0.1356091046061895
$ws.Range("AB7").Value = -603.3432340966085
$ws.Range("AC7").Value = "umolO2/min/m2"
$ws.Range("AD7").Value = -603.3432340966085

$ws.Range("T8").Value = 0.1544
$ws.Range("V8").Value = 0
$ws.Range("Z8").Value = 0.007024329983665378
$ws.Range("AB8").Value = "Inf"
$ws.Range("AC8").Value = "umolO2/min/m2"
$ws.Range("AD8").Value = "Inf"

$ws.Range("T9").Value = 0.1426048780487805
$ws.Range("V9").Value = 0.0002220611491829204
$ws.Range("Z9").Value = 0.02794665006595875
$ws.Range("AB9").Value = 125.8511458163175
$ws.Range("AC9").Value = "umolO2/min/m2"
$ws.Range("AD9").Value = 125.8511458163175

$ws.Range("T10").Value = 0.1423414634146342
$ws.Range("V10").Value = 0.0002565454225970831
$ws.Range("Z10").Value = -0.001984993206628306
$ws.Range("AB10").Value = -7.737394752686086
$ws.Range("AC10").Value = "umolO2/min/m2"
$ws.Range("AD10").Value = -7.737394752686086

$ws.Range("T11").Value = 0.1477268292682927
$ws.Range("V11").Value = 0.0002491214197856265
$ws.Range("Z11").Value = 0.1211464226602152
$ws.Range("AB11").Value = 486.29468619946
$ws.Range("AC11").Value = "umolO2/min/m2"
$ws.Range("AD11").Value = 486.29468619946

$ws.Range("T12").Value = 0.1470341463414634
$ws.Range("V12").Value = 0.0002082674398172554
$ws.Range("Z12").Value = 0.1044969838570853
$ws.Range("AB12").Value = 501.7442186295485
$ws.Range("AC12").Value = "umolO2/min/m2"
$ws.Range("AD12").Value = 501.7442186295485

$ws.Range("T13").Value = 0.1477560975609756
$ws.Range("V13").Value = 0.0001607142857142857
$ws.Range("Z13").Value = 0.1102646387779767
$ws.Range("AB13").Value = 686.0910857296326
$ws.Range("AC13").Value = "umolO2/min/m2"
$ws.Range("AD13").Value = 686.0910857296326

$ws.Range("T14").Value = 0.1455317073170732
$ws.Range("V14").Value = 0.0002247627833421192
$ws.Range("Z14").Value = -0.02142235960129866
$ws.Range("AB14").Value = -95.31097311911708
$ws.Range("AC14").Value = "umolO2/min/m2"
$ws.Range("AD14").Value = -95.31097311911708

$ws.Range("T15").Value = 0.1544
$ws.Range("V15").Value = 0
$ws.Range("Z15").Value = -0.0006681919434117178
$ws.Range("AB15").Value = "-Inf"
$ws.Range("AC15").Value = "umolO2/min/m2"
$ws.Range("AD15").Value = "-Inf"
